$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (after the edit) for rows 2-51, columns A:C
$data = @(
    @(299.9468443580553, 1.9886474699379992, -0.9636595695610491),
    @(311.14773847666396, 2.0049755339883277, -1.1267831787260159),
    @(293.5900191095634, 2.21276954675977, -0.988400062491897),
    @(305.4400252289269, 1.902956612567596, -1.0405615097602228),
    @(302.7973681419485, 1.9350478927000663, -0.9897465263559363),
    @(278.79106516956915, 2.265801802008292, -0.9295754784372078),
    @(314.6212340033373, 1.8928606502430152, -1.2104914292014124),
    @(288.8306694108422, 1.9930073668326012, -0.881154869614634),
    @(293.28128173565705, 2.016476273448693, -0.8820076306874003),
    @(314.9250197665565, 1.887468839036188, -1.1800383197781301),
    @(297.1174698980124, 1.7657724666533703, -0.7987880547195967),
    @(291.09854555879275, 1.9984892611441916, -0.839598138763102),
    @(303.700650965937, 1.7346217626892675, -0.9575360054727229),
    @(296.4304785032986, 1.8421916992631149, -0.9056114803187651),
    @(295.53522825552216, 2.0270269595536545, -1.005249690365148),
    @(301.066060961494, 1.9811191871485336, -0.9583768516208291),
    @(298.99288924825606, 2.0799064181379543, -0.9671965434765146),
    @(289.40986701228513, 2.0821628565018333, -0.9730753862045277),
    @(295.8612060366306, 1.9638773969950263, -1.0201240474488815),
    @(308.14191260070356, 1.9502000370830177, -1.0845055422170364),
    @(312.36158519753366, 1.82896550853489, -1.0330795358159008),
    @(309.02331515811625, 1.5039832877679, -0.8369411233383217),
    @(301.7003658475296, 1.9131902769816842, -0.9987340104937925),
    @(306.3086131837516, 1.933575172159937, -1.1207670675435606),
    @(311.1954306212756, 1.829971141213248, -1.0513841715027088),
    @(299.72951955610176, 1.9264166811862933, -0.7951797454972417),
    @(308.0242603359436, 1.8484149705322845, -1.061275713332102),
    @(292.5140104951312, 1.8840618460744176, -0.7410037380078798),
    @(316.50602659791673, 1.7444387271296782, -1.0862605418301323),
    @(306.0487930540122, 1.8798858534006535, -1.0607872934002338),
    @(281.1177061105754, 2.052272695254657, -0.7048513578789039),
    @(310.4497019125892, 1.9283381295106592, -1.0443681163572494),
    @(310.84432884340345, 1.9604232119222391, -1.1839624811767226),
    @(289.91348699104594, 2.135380586900032, -1.0159315718831692),
    @(291.9735544212472, 2.1560368953440032, -1.011364365655876),
    @(295.51420590818935, 2.1481206961008024, -1.0601778303136262),
    @(287.65985692220744, 2.0760010753148412, -0.8808939758662092),
    @(304.79924631630905, 1.8160360778679385, -0.8078139382270659),
    @(303.6912216786658, 2.033844966426333, -1.1259242684460369),
    @(298.7671609537739, 2.206863718022714, -1.1194197679033593),
    @(300.7885664154631, 1.8682658018213865, -0.9077520023487906),
    @(295.7682577602744, 2.044769143187393, -0.9708087664457066),
    @(289.57519440403996, 2.1824352648702066, -0.9582761577790398),
    @(296.018820746225, 2.056410150646353, -0.9148865468657899),
    @(292.16458531584686, 2.023578355054992, -1.0284014895437639),
    @(306.00736653870547, 2.051365703630159, -1.213766290622896),
    @(308.19925715442076, 2.034659915275052, -1.0828098329640194),
    @(296.90937478287304, 1.9222067547404746, -0.8530735465760723),
    @(305.44083406794346, 2.0235087111515493, -1.1304663743678722),
    @(313.1010051940947, 1.8828182420302344, -1.0205707631826801)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Original sheet had 52 data rows (row 2-52); new layout only has 50 data rows (2-51).
# Remove the now-obsolete last row (52) entirely.
$ws.Range("A52:C52").EntireRow.Delete() | Out-Null
